# Weekly update: insert two new price-report rows (this week's Mango entries)
# right after the current top entries, pushing the rest of the history down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at row 76; existing rows 76-117 shift down to 78-119.
$ws.Rows.Item(76).Resize(2).Insert()

# New row 76: Primera quality entry for the new week.
$ws.Cells.Item(76, 1).Value  = 4
$ws.Cells.Item(76, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(76, 3).Value  = "Los Lagos"
$ws.Cells.Item(76, 4).Value  = [datetime]"2021-11-29"
$ws.Cells.Item(76, 5).Value  = 10
$ws.Cells.Item(76, 6).Value  = "Fruta"
$ws.Cells.Item(76, 7).Value  = 100108
$ws.Cells.Item(76, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(76, 9).Value  = 100108002
$ws.Cells.Item(76, 10).Value = "Mango"
$ws.Cells.Item(76, 11).Value = "Sin especificar"
$ws.Cells.Item(76, 12).Value = "Primera"
$ws.Cells.Item(76, 13).Value = 60
$ws.Cells.Item(76, 14).Value = 7500
$ws.Cells.Item(76, 15).Value = 8000
$ws.Cells.Item(76, 16).Value = 7750
$ws.Cells.Item(76, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(76, 18).Value = "Perú"
$ws.Cells.Item(76, 19).Value = 1938
$ws.Cells.Item(76, 20).Value = 4

# New row 77: Segunda quality entry for the new week.
$ws.Cells.Item(77, 1).Value  = 4
$ws.Cells.Item(77, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(77, 3).Value  = "Los Lagos"
$ws.Cells.Item(77, 4).Value  = [datetime]"2021-11-29"
$ws.Cells.Item(77, 5).Value  = 10
$ws.Cells.Item(77, 6).Value  = "Fruta"
$ws.Cells.Item(77, 7).Value  = 100108
$ws.Cells.Item(77, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(77, 9).Value  = 100108002
$ws.Cells.Item(77, 10).Value = "Mango"
$ws.Cells.Item(77, 11).Value = "Sin especificar"
$ws.Cells.Item(77, 12).Value = "Segunda"
$ws.Cells.Item(77, 13).Value = 40
$ws.Cells.Item(77, 14).Value = 5000
$ws.Cells.Item(77, 15).Value = 5000
$ws.Cells.Item(77, 16).Value = 5000
$ws.Cells.Item(77, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(77, 18).Value = "Perú"
$ws.Cells.Item(77, 19).Value = 1250
$ws.Cells.Item(77, 20).Value = 4
